$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ A="Inflammatory-Mac"; B="Ifng"; C="Ifngr2"; D="ECs"; E=2; F=0.6666666666666666; G=0.4227963333333333; H=1.268389; I=0.9187469758736285; J=0.9187469758736285; K=3; L=1; M=5.595533666666667; N=16.786601; O=0.04791382214400042; P=0.048561535774258; Q=2.365771117309889; R=21.291940055789; S=0.04402067919734728; T=0.04461576413637856 },
  @{ A="Inflammatory-Mac"; B="Ifng"; C="Ifngr2"; D="FAPs"; E=2; F=0.6666666666666666; G=0.4227963333333333; H=1.268389; I=0.9187469758736285; J=0.9187469758736285; K=3; L=1; M=10.82716733333333; N=32.481502; O=0.09271161623475734; P=0.09396492007968932; Q=4.577686648919778; R=41.199179840278; S=0.0851785170440397; T=0.08632998616142175 },
  @{ A="Inflammatory-Mac"; B="Ifng"; C="Ifngr2"; D="Inflammatory-Mac"; E=2; F=0.6666666666666666; G=0.4227963333333333; H=1.268389; I=0.9187469758736285; J=0.9187469758736285; K=3; L=1; M=52.932222; N=158.796666; O=0.4532516863767858; P=0.4593788806198408; Q=22.379549376786; R=201.415944391074; S=0.4164236161682942; T=0.4220529573496913 },
  @{ A="Inflammatory-Mac"; B="Ifng"; C="Ifngr2"; D="MuSCs"; E=2; F=0.6666666666666666; G=0.4227963333333333; H=1.268389; I=0.9187469758736285; J=0.9187469758736285; K=2; L=1; M=4.6729655; N=9.345931; O=0.04001399172805441; P=0.02703660869763014; Q=1.975712679193167; R=11.854276075159; S=0.03676273389278238; T=0.02483980247882633 },
  @{ A="Inflammatory-Mac"; B="Ifng"; C="Ifngr2"; D="Resolving-Mac"; E=2; F=0.6666666666666666; G=0.4227963333333333; H=1.268389; I=0.9187469758736285; J=0.9187469758736285; K=3; L=1; M=42.755399; N=128.266197; O=0.3661088835164021; P=0.3710580548285817; Q=18.07682592740367; R=162.691433346633; S=0.336361429571165; T=0.3409084657473104 },
  @{ A="Resolving-Mac"; B="Ifng"; C="Ifngr2"; D="ECs"; E=1; F=0.3333333333333333; G=0.03739166666666666; H=0.112175; I=0.08125302412637153; J=0.08125302412637155; K=3; L=1; M=5.595533666666667; N=16.786601; O=0.04791382214400042; P=0.048561535774258; Q=0.2092263296861111; R=1.883036967175; S=0.00389314294665314; T=0.00394577163787944 },
  @{ A="Resolving-Mac"; B="Ifng"; C="Ifngr2"; D="FAPs"; E=1; F=0.3333333333333333; G=0.03739166666666666; H=0.112175; I=0.08125302412637153; J=0.08125302412637155; K=3; L=1; M=10.82716733333333; N=32.481502; O=0.09271161623475734; P=0.09396492007968932; Q=0.4048458318722222; R=3.64361248685; S=0.007533099190717637; T=0.00763493391826757 },
  @{ A="Resolving-Mac"; B="Ifng"; C="Ifngr2"; D="Inflammatory-Mac"; E=1; F=0.3333333333333333; G=0.03739166666666666; H=0.112175; I=0.08125302412637153; J=0.08125302412637155; K=3; L=1; M=52.932222; N=158.796666; O=0.4532516863767858; P=0.4593788806198408; Q=1.97922400095; R=17.81301600855; S=0.03682807020849155; T=0.03732592327014948 },
  @{ A="Resolving-Mac"; B="Ifng"; C="Ifngr2"; D="MuSCs"; E=1; F=0.3333333333333333; G=0.03739166666666666; H=0.112175; I=0.08125302412637153; J=0.08125302412637155; K=2; L=1; M=4.6729655; N=9.345931; O=0.04001399172805441; P=0.02703660869763014; Q=0.1747299683208333; R=1.048379809925; S=0.003251257835272036; T=0.002196806218803809 },
  @{ A="Resolving-Mac"; B="Ifng"; C="Ifngr2"; D="Resolving-Mac"; E=1; F=0.3333333333333333; G=0.03739166666666666; H=0.112175; I=0.08125302412637153; J=0.08125302412637155; K=3; L=1; M=42.755399; N=128.266197; O=0.3661088835164021; P=0.3710580548285817; Q=1.598695627608333; R=14.388260648475; S=0.02974745394523717; T=0.03014958908127124 }
)

$cols = @{ A=1; B=2; C=3; D=4; E=5; F=6; G=7; H=8; I=9; J=10; K=11; L=12; M=13; N=14; O=15; P=16; Q=17; R=18; S=19; T=20 }
$colLetters = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$r = 2
foreach ($row in $rows) {
    foreach ($col in $colLetters) {
        $ws.Cells.Item($r, $cols[$col]).Value = $row[$col]
    }
    $r++
}
